$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.872.53"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.389.29"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "580.19"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "178.82"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +4.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.387.78"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "3.979.50"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "28.89"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").Value = "66.013.04"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "3.391.15"
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("D19").Value = "5.85"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "13.67"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "365.09"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "7.50"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").Value = "72.78"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "5.70"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "23.05"
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "6.95"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "161.11"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").Value = "0.859"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").Value = "27.13"
$ws.Range("E39").Value = "  -8.05%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "2.664.24"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "6.20"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "0.0676"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "39.59"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "24.32"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "331.38"
$ws.Range("E48").Value = "  +7.71%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("D51").Value = "31.32"
$ws.Range("E51").Value = "  +3.65%  "
